# Update cryptocurrency price/volume data per latest GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.415.81'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.725.61'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.80'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4900'
$ws.Range('E7').Value = '  +1.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2608'
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06198'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.722.71'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07017'
$ws.Range('E11').Value = '  -2.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.53'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6004'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.35'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.431.31'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007147'
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.35'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.942.55'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.474'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.591'
$ws.Range('E23').Value = '  -2.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.167'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.68'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.23'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.392'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.92'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.703'
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.960'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07958'
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.684'
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04533'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.604'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9973'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6275'
$ws.Range('E36').Value = '  +0.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9118'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.960'
$ws.Range('E38').Value = '  -5.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.393'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.43'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.445'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3849'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.719'
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1159'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05367'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.738'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.14'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.239'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.06'
$ws.Range('E51').Value = '  -0.13%  '
